$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "Wins", "Losses", "Ties" columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold/border/centered) from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill team record values for all data rows (2 through 43)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 87
    $ws.Cells.Item($r, 31).Value = 75
    $ws.Cells.Item($r, 32).Value = 0
}
